$d = $word.ActiveDocument

# The edit only touches the first paragraph of the document (the one holding the
# **ID__AFFARS_MP5301_90__ID** merge placeholder):
#   1. A paragraph border is added, with each side (top/left/bottom/right) set to
#      a 5pt "space" (distance from text) but no visible line.
#   2. The paragraph's left indent changes from 120 twips (6pt) to 225 twips (11.25pt).
#   3. The trailing run that contained just a single space " " is removed, leaving
#      only the placeholder run.

$p1 = $d.Paragraphs.Item(1)

# 2. Update left indent (225 twips = 11.25 points).
$p1.Format.LeftIndent = 11.25

# 1. Add paragraph border spacing (5 points on every side).
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# 3. Remove the trailing " " run by replacing the paragraph's text (placeholder + space)
#    with just the placeholder text; this merges away the now-empty second run.
$rng = $p1.Range
$rng.Find.Execute("**ID__AFFARS_MP5301_90__ID** ", $true, $false, $false, $false, $false,
                   $true, 1, $false, "**ID__AFFARS_MP5301_90__ID**", 2)

Write-Host "Paragraph 1 updated: border + indent + trailing run removed"
